# Update RMIL data 2025-11-01 03:55:45 UTC
#
# 1) Four cells in the existing "2025-10-27" block were re-typed from
#    inline-string "0" to a real numeric 0 (H56, H57, F58, H58).
# 2) A new circular dated 2025-10-30 was appended as rows 70-86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-type four text "0" cells as numeric 0 -----------------------
$ws.Range("H56").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("H58").Value = 0

# --- 2) Append new circular rows 70-86 ----------------------------------
# Force text format on the new data range so numeric-looking strings
# (dates, comma-grouped numbers, "0") are stored as text, matching the source data.
$ws.Range("A70:I86").NumberFormat = "@"

# Row 70
$ws.Range("A70").Value = "2025-10-30"
$ws.Range("B70").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C70").Value = "CHANDERIYA `nLEAD ZINC `nSMELTER"
$ws.Range("E70").Value = "330,100 331,600 330,600 329,600 328,100"
$ws.Range("I70").Value = "209,400"

# Row 71
$ws.Range("A71").Value = "2025-10-30"
$ws.Range("B71").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C71").Value = "HYDRO-1 UNIT"
$ws.Range("E71").Value = "330,100 331,600 330,600 329,600 328,100"
$ws.Range("I71").Value = "209,400"

# Row 72
$ws.Range("A72").Value = "2025-10-30"
$ws.Range("B72").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C72").Value = "NEW HYDRO `nSMELTER `nCHANDERIYA"
$ws.Range("E72").Value = "330,100 331,600 330,600 329,600 328,100"
$ws.Range("I72").Value = "209,400"

# Row 73
$ws.Range("A73").Value = "2025-10-30"
$ws.Range("B73").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C73").Value = "ZINC SMELTER `nDEBRI"
$ws.Range("D73").Value = "0"
$ws.Range("E73").Value = "0"
$ws.Range("G73").Value = "0  329,600"
$ws.Range("H73").Value = "0"
$ws.Range("I73").Value = "0"

# Row 74
$ws.Range("A74").Value = "2025-10-30"
$ws.Range("B74").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C74").Value = "Pantnagar `nMelting&Castin `ngPlant"
$ws.Range("E74").Value = "330,100 331,600 330,600 329,600"
$ws.Range("H74").Value = "0"
$ws.Range("I74").Value = "209,400"

# Row 75
$ws.Range("A75").Value = "2025-10-30"
$ws.Range("B75").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C75").Value = "RAJPURA DARIBA `nLEAD SMELTER"
$ws.Range("D75").Value = "0"
$ws.Range("E75").Value = "0"
$ws.Range("F75").Value = "0"
$ws.Range("G75").Value = "0"
$ws.Range("H75").Value = "0"
$ws.Range("I75").Value = "209,400"

# Row 76
$ws.Range("A76").Value = "2025-10-30"
$ws.Range("B76").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C76").Value = "Faridabad `nDepot"
$ws.Range("E76").Value = "332,600 334,100 328,100 332,100 330,600"
$ws.Range("I76").Value = "211,900"

# Row 77
$ws.Range("A77").Value = "2025-10-30"
$ws.Range("B77").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C77").Value = "Panvel Depot"
$ws.Range("E77").Value = "333,400 334,900 333,900 332,900 331,400"
$ws.Range("I77").Value = "212,300"

# Row 78
$ws.Range("A78").Value = "2025-10-30"
$ws.Range("B78").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C78").Value = "Pune Depot"
$ws.Range("E78").Value = "333,400 334,900 333,900 332,900 331,400"
$ws.Range("I78").Value = "212,700"

# Row 79
$ws.Range("A79").Value = "2025-10-30"
$ws.Range("B79").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C79").Value = "Baroda Depot"
$ws.Range("E79").Value = "333,400 334,900 333,900 332,900 331,400"
$ws.Range("I79").Value = "212,700"

# Row 80
$ws.Range("A80").Value = "2025-10-30"
$ws.Range("B80").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C80").Value = "Raipur Depot"
$ws.Range("E80").Value = "333,400 334,900 333,900 332,900 331,400"
$ws.Range("I80").Value = "212,700"

# Row 81
$ws.Range("A81").Value = "2025-10-30"
$ws.Range("B81").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C81").Value = "JAMSHEDPUR `nSTOCK POINT"
$ws.Range("E81").Value = "331,100 332,600 331,600 330,600 329,100"
$ws.Range("I81").Value = "210,400"

# Row 82
$ws.Range("A82").Value = "2025-10-30"
$ws.Range("B82").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("D82").Value = "Kolkata Depot  331,100 332,600 331,600 330,600 329,100"
$ws.Range("I82").Value = "210,400"

# Row 83
$ws.Range("A83").Value = "2025-10-30"
$ws.Range("B83").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C83").Value = "Bangalore `nDepot"
$ws.Range("E83").Value = "331,100 332,600 331,600 330,600 329,100"
$ws.Range("I83").Value = "210,400"

# Row 84
$ws.Range("A84").Value = "2025-10-30"
$ws.Range("B84").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C84").Value = "Hyderabad `nDepot"
$ws.Range("E84").Value = "331,100 332,600 331,600 330,600 329,100"
$ws.Range("I84").Value = "210,400"

# Row 85
$ws.Range("A85").Value = "2025-10-30"
$ws.Range("B85").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("D85").Value = "Chennai Depot  331,100 332,600 331,600 330,600 329,100"
$ws.Range("I85").Value = "210,400"

# Row 86
$ws.Range("A86").Value = "2025-10-30"
$ws.Range("B86").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C86").Value = "Sindesar `nsmelter HZAPL"
$ws.Range("E86").Value = "330,100 331,600"
$ws.Range("G86").Value = "0  329,600 328,100"
$ws.Range("I86").Value = "209,400"
